$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Build C45 first (font6 black Calibri/12)
$ws.Range("C44").Copy()
$ws.Range("C45").PasteSpecial(-4122)
$ws.Range("C45").Value = "CopyCat"
$ws.Range("C45").Font.Color = 0
Write-Output "C45: font=$($ws.Range('C45').Font.Name)/$($ws.Range('C45').Font.Size)/$($ws.Range('C45').Font.Color) fmt=$($ws.Range('C45').NumberFormat)"

# Now build D45 by copying the STYLE from C45 (which has font6, numFmt0), then set numberformat/value
$ws.Range("C45").Copy()
$ws.Range("D45").PasteSpecial(-4122)
Write-Output "D45 after format paste from C45: font=$($ws.Range('D45').Font.Name)/$($ws.Range('D45').Font.Size)/$($ws.Range('D45').Font.Color) fmt=$($ws.Range('D45').NumberFormat)"
$ws.Range("D45").Value2 = 42879
$ws.Range("D45").NumberFormat = "m/d/yyyy"
Write-Output "D45 final: font=$($ws.Range('D45').Font.Name)/$($ws.Range('D45').Font.Size)/$($ws.Range('D45').Font.Color) fmt=$($ws.Range('D45').NumberFormat)"
